# Add "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header formatting from the neighboring "sum" header cell (G1)
# so the new "Save" header (H1) picks up the same style index instead of
# Excel minting a brand-new style/font combination.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data column values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
